$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet "2022-Q4" right after "总计" (before the
#        existing "2021-Q4" sheet), mirroring the layout of the existing
#        per-quarter fund sheet. ---
$totalSheet = $wb.Worksheets.Item(1)
$oldQuarterSheet = $wb.Worksheets.Item(2)

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q4"

# Header row (same headers as the existing "2021-Q4" fund sheet)
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$hdr = $newSheet.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data row
$a2 = $newSheet.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "159628"
$newSheet.Range("B2").ClearFormats()

$newSheet.Range("C2").Value = "万家国证2000ETF"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.55"
$newSheet.Range("D2").ClearFormats()

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "97.28"
$newSheet.Range("E2").ClearFormats()

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "0.45"
$newSheet.Range("F2").ClearFormats()

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0115"
$newSheet.Range("G2").ClearFormats()

$newSheet.Range("H2").Value = 7

# --- 2. "总计" summary sheet: the row that used to describe 2021-Q4 now
#        describes 2022-Q4, and a new row is appended for 2021-Q4 (pushed
#        down to row 3). ---
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A3").Font.Bold = $true
$totalSheet.Range("A3").HorizontalAlignment = -4108
$totalSheet.Range("A3").VerticalAlignment = -4160
$totalSheet.Range("A3").Borders.LineStyle = 1

$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.15
